$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in I1 from "eta_tilde_engr_units" to "eta_star_engr_units"
$ws.Range("I1").Value = "eta_star_engr_units"

# Update the selected cell to I1 (matches the new active selection in the diff)
$ws.Range("I1").Select()
